# Day 1 contents update
# 1) Slide 6 ("Day 4" bullet): split the run that contains
#    ": ML models, from its selection to its training, evaluation and tuning"
#    into two runs - "...evaluation and " / "tuning" - with identical formatting.
# 2) Theme color refresh (Design colors changed from "Simple Light" to the
#    "Default" palette on the presentation's addressable theme).

$p = $ppt.ActivePresentation

# ---- 1) Split the "tuning" run on the Day 4 bullet (slide 6) ----
$slide = $p.Slides.Item(6)

$introShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $candidate = $slide.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText) {
        if ($candidate.TextFrame.TextRange.Text -like "*ML models, from its selection*") {
            $introShape = $candidate
        }
    }
}

$tr = $introShape.TextFrame.TextRange

$targetPara = $null
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text -like "*ML models, from its selection*") {
        $targetPara = $para
    }
}

$paraText = $targetPara.Text
$word = "tuning"
$wordOffset = $paraText.LastIndexOf($word)

$wordStart = $targetPara.Start + $wordOffset
$wordRange = $tr.Characters($wordStart, $word.Length)
# Re-assigning the text of this sub-range (identical text, same formatting)
# forces PowerPoint to materialize it as its own run, splitting the
# previously single run in two while keeping matching run properties.
$wordRange.Text = $word

# ---- 2) Theme color palette update ----
$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
$newPalette = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    5800213,    # dk2      158158
    15987699,   # lt2      F3F3F3
    13077765,   # accent1  058DC7
    3322960,    # accent2  50B432
    1791725,    # accent3  ED561B
    61421,      # accent4  EDEF00
    15059748,   # accent5  24CBE5
    7529828,    # accent6  64E572
    13369378,   # hlink    2200CC
    9116245     # folHlink 551A8B
)
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $newPalette[$i - 1]
}
